# Update the Gantt chart data/view as described in the commit:
# - B8 (Open) 39 -> 38
# - C8 (Closed) 12 -> 13
# - Selection moves from D8 to G10
# - Workbook window width grows (window resized wider)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Board")

$ws.Range("B8").Value = 38
$ws.Range("C8").Value = 13

$ws.Activate()
$ws.Range("G10").Select()

$excel.ActiveWindow.Width = 33120
